# Actualización automática 2025-12-10 08:30:06
# Applies updated sales figures to "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("M4").Value = 1287.85
$wsGrupo.Range("D5").Value = 933.12
$wsGrupo.Range("D7").Value = 380.16
$wsGrupo.Range("C12").Value = 388.8
$wsGrupo.Range("K12").Value = 1268.64
$wsGrupo.Range("L12").Value = 1393.76
$wsGrupo.Range("M19").Value = 1833.06
$wsGrupo.Range("K36").Value = 92.88
$wsGrupo.Range("M36").Value = 9531.92
$wsGrupo.Range("P36").Value = 489.75
$wsGrupo.Range("M37").Value = 1960.7
$wsGrupo.Range("M53").Value = 917.8

# Row 56 contains the "X de 54" completion counters.
$wsGrupo.Range("C56").Value = "2 de 54"
$wsGrupo.Range("D56").Value = "4 de 54"
$wsGrupo.Range("K56").Value = "2 de 54"
$wsGrupo.Range("M56").Value = "11 de 54"

# --- Sheet: VENTA MENSUAL ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F4").Value = 1642.09
$wsMensual.Range("F5").Value = 2476.49
$wsMensual.Range("F7").Value = 380.16
$wsMensual.Range("F12").Value = 3051.2
$wsMensual.Range("F19").Value = 2973.38
$wsMensual.Range("F36").Value = 10657.25
$wsMensual.Range("F37").Value = 4128.8
$wsMensual.Range("F55").Value = 1207.6
$wsMensual.Range("F56").Value = 1207.6
$wsMensual.Range("F60").Value = 34426.16
